$d = $word.ActiveDocument

# Find the paragraph holding the final bibliography entry ("...Thomson,
# 2007."). The three paragraphs immediately following it are the site
# "footer" block that needs to go: a blank spacer paragraph, the
# "Ver no Jupiter..." line, and the "(c) 2020 ..." copyright line.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Thomson, 2007.*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $firstToDelete = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToDelete = $d.Paragraphs.Item($anchorIndex + 3)

    $deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
    $deleteRange.Delete()
}
